$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.501.73"
$ws.Range("E2").Value = "  +4.53%  "
$ws.Range("D3").Value = "1.735.68"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4791"
$ws.Range("E7").Value = "  +3.95%  "
$ws.Range("D8").Value = "0.2662"
$ws.Range("E8").Value = "  +3.98%  "
$ws.Range("D9").Value = "0.06224"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "1.736.05"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "15.71"
$ws.Range("E12").Value = "  +8.17%  "
$ws.Range("D13").Value = "0.6126"
$ws.Range("E13").Value = "  +7.56%  "
$ws.Range("D14").Value = "4.525"
$ws.Range("E14").Value = "  +5.03%  "
$ws.Range("D15").Value = "76.82"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "26.513.21"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "0.000006897"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").Value = "1.959.94"
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("D22").Value = "4.565"
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("D23").Value = "8.878"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "5.331"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "15.35"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").Value = "1.795"
$ws.Range("E27").Value = "  +5.64%  "
$ws.Range("D28").Value = "1.401"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "106.82"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "3.977"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.07948"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.706"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("D33").Value = "0.04555"
$ws.Range("E33").Value = "  +5.51%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6359"
$ws.Range("E36").Value = "  +6.52%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.9908"
$ws.Range("E37").Value = "  +6.38%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.9340"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "110.71"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.983"
$ws.Range("E40").Value = "  +9.30%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.416"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01508"
$ws.Range("E43").Value = "  +4.30%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.692"
$ws.Range("E44").Value = "  +14.66%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3895"
$ws.Range("E45").Value = "  +5.80%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.923"
$ws.Range("E46").Value = "  +14.13%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1190"
$ws.Range("E47").Value = "  +8.34%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05335"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.910"
$ws.Range("E49").Value = "  +4.91%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.76"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.255"
$ws.Range("E51").Value = "  +6.13%  "
